# Apply updated cryptocurrency price/volume data to match the
# Sun Aug 13 05:37:03 UTC 2023 GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.405.04'
$ws.Range("D3").Value = '''1.849.78'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D5").Value = '''240.51'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = '''0.6273'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").Value = '''0.07636'
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("D9").Value = '''0.2906'
$ws.Range("E9").Value = '  -0.78%  '
$ws.Range("D10").Value = '''24.76'
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("D11").Value = '''0.07749'
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").Value = '''5.036'
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("D13").Value = '''0.6790'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '''0.00001066'
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("D15").Value = '''83.25'
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = '''6.159'
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = '''29.423.34'
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("D18").Value = '''226.44'
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("E19").Value = '  -0.89%  '
$ws.Range("D20").Value = '''1.002'
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").Value = '''7.473'
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '''157.93'
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").Value = '''0.1379'
$ws.Range("E24").Value = '  -1.11%  '
$ws.Range("D25").Value = '''8.418'
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").Value = '''17.68'
$ws.Range("D27").Value = '''1.387'
$ws.Range("E27").Value = '  +6.69%  '
$ws.Range("D28").Value = '''1.464'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '''0.05587'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").Value = '''4.127'
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("D31").Value = '''4.056'
$ws.Range("E31").Value = '  +0.55%  '
$ws.Range("D32").Value = '''1.837'
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("D33").Value = '''1.163'
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("D34").Value = '''0.6946'
$ws.Range("E34").Value = '  -1.95%  '
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '''1.233.11'
$ws.Range("E36").Value = '  +0.29%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.01803'
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("D38").Value = '''2.720'
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("D39").Value = '''6.414'
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("D40").Value = '''0.9040'
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D42").Value = '''101.68'
$ws.Range("E42").Value = '  -0.20%  '
$ws.Range("D43").Value = '''65.95'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value = '''0.00000000121'
$ws.Range("E44").Value = '  -0.90%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '''7.173'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '''0.4011'
$ws.Range("E46").Value = '  -0.14%  '
$ws.Range("D47").Value = '''8.976'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").Value = '''0.1144'
$ws.Range("E49").Value = '  +2.05%  '
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("E51").Value = '  +0.23%  '
